$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$startDevice = 3000166

for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $startDevice + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

$ws.Range("C152").Select()
$excel.ActiveWindow.ScrollRow = 140
